# "Generate Report for Handback"
#
# The localization report's "3de4c1a0-...md" row has moved from
# "awaiting handoff" to "handed back" (in sync with en-US):
#   - Overview sheet: zh-cn / de-de status columns flip to the
#     "Handed back: in sync with en-US" status text.
#   - zh-cn / de-de detail sheets: Status flips the same way, the
#     "Latest Handback DateTime" is refreshed to the new handback
#     timestamp, and the stale "Error Detail" message is cleared
#     now that the handback succeeded.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---- zh-cn detail sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-09-03 02:52:52"
$wsZhCn.Range("P3").ClearContents()
$wsZhCn.Columns.Item(16).ColumnWidth = 12.86

# ---- de-de detail sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-09-03 02:52:59"
$wsDeDe.Range("P3").ClearContents()
$wsDeDe.Columns.Item(16).ColumnWidth = 12.86
